$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 15, shifting existing rows 15-23 down to 16-24
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new data record
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44452
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 100112026
$ws.Range("G15").Value = "Haba"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 13000
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 520
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
